$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 3265880  # ALC!H132 (4489999 -> 3265880)
$ws.Cells.Item(132, 9).Value = 3592168.8  # ALC!I132 (4789099.5 -> 3592168.8)
$ws.Cells.Item(132, 10).Value = 2992.75  # ALC!J132 (3498.5 -> 2992.75)
$ws.Cells.Item(132, 11).Value = 10776506.4  # ALC!K132 (14367298.5 -> 10776506.4)
$ws.Cells.Item(132, 12).Value = 8978.25  # ALC!L132 (10495.5 -> 8978.25)
$ws.Cells.Item(132, 13).Value = -10773976.4  # ALC!M132 (-14364768.5 -> -10773976.4)
$ws.Cells.Item(132, 14).Value = -14038.25  # ALC!N132 (-15555.5 -> -14038.25)

$ws.Cells.Item(135, 8).Value = 13482.174  # ALC!H135 (59950.4 -> 13482.174)
$ws.Cells.Item(135, 9).Value = 753.6  # ALC!I135 (2622.6667 -> 753.6)
$ws.Cells.Item(135, 10).Value = 98339.336  # ALC!J135 (84519.42999999999 -> 98339.336)
$ws.Cells.Item(135, 11).Value = 6782.400000000001  # ALC!K135 (23604.0003 -> 6782.400000000001)
$ws.Cells.Item(135, 12).Value = 885054.024  # ALC!L135 (760674.8699999999 -> 885054.024)
$ws.Cells.Item(135, 13).Value = -4247.400000000001  # ALC!M135 (-21069.0003 -> -4247.400000000001)
$ws.Cells.Item(135, 14).Value = -890124.024  # ALC!N135 (-765744.8699999999 -> -890124.024)

$ws.Cells.Item(137, 8).Value = 13109.241  # ALC!H137 (15023.12 -> 13109.241)
$ws.Cells.Item(137, 9).Value = 17940.63  # ALC!I137 (23805.643 -> 17940.63)
$ws.Cells.Item(137, 10).Value = 3929.6  # ALC!J137 (3845.3635 -> 3929.6)
$ws.Cells.Item(137, 11).Value = 53821.89  # ALC!K137 (71416.929 -> 53821.89)
$ws.Cells.Item(137, 12).Value = 11788.8  # ALC!L137 (11536.0905 -> 11788.8)
$ws.Cells.Item(137, 13).Value = -51271.89  # ALC!M137 (-68866.929 -> -51271.89)
$ws.Cells.Item(137, 14).Value = -16888.8  # ALC!N137 (-16636.0905 -> -16888.8)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 18895.908  # ARM!H32 (18962.969 -> 18895.908)
$ws.Cells.Item(32, 9).Value = 19508.412  # ARM!I32 (19588.598 -> 19508.412)
$ws.Cells.Item(32, 11).Value = 19508.412  # ARM!K32 (19588.598 -> 19508.412)
$ws.Cells.Item(32, 13).Value = -19221.412  # ARM!M32 (-19301.598 -> -19221.412)

$ws.Cells.Item(61, 8).Value = 6713.0415  # ARM!H61 (4456.4873 -> 6713.0415)
$ws.Cells.Item(61, 9).Value = 1711.2106  # ARM!I61 (1203.1212 -> 1711.2106)
$ws.Cells.Item(61, 10).Value = 25720  # ARM!J61 (22350 -> 25720)
$ws.Cells.Item(61, 11).Value = 1711.2106  # ARM!K61 (1203.1212 -> 1711.2106)
$ws.Cells.Item(61, 12).Value = 25720  # ARM!L61 (22350 -> 25720)
$ws.Cells.Item(61, 13).Value = -1499.2106  # ARM!M61 (-991.1212 -> -1499.2106)
$ws.Cells.Item(61, 14).Value = -26144  # ARM!N61 (-22774 -> -26144)

$ws.Cells.Item(74, 8).Value = 125428.266  # ARM!H74 (139652.16 -> 125428.266)
$ws.Cells.Item(74, 9).Value = 125958.02  # ARM!I74 (140574.3 -> 125958.02)
$ws.Cells.Item(74, 11).Value = 125958.02  # ARM!K74 (140574.3 -> 125958.02)
$ws.Cells.Item(74, 13).Value = -125084.02  # ARM!M74 (-139700.3 -> -125084.02)

$ws.Cells.Item(77, 8).Value = 125428.266  # ARM!H77 (139652.16 -> 125428.266)
$ws.Cells.Item(77, 9).Value = 125958.02  # ARM!I77 (140574.3 -> 125958.02)
$ws.Cells.Item(77, 11).Value = 629790.1  # ARM!K77 (702871.5 -> 629790.1)
$ws.Cells.Item(77, 13).Value = -625422.1  # ARM!M77 (-698503.5 -> -625422.1)

$ws.Cells.Item(122, 8).Value = 1975.9487  # ARM!H122 (1699.1136 -> 1975.9487)
$ws.Cells.Item(122, 9).Value = 1932.2778  # ARM!I122 (1692.5128 -> 1932.2778)
$ws.Cells.Item(122, 10).Value = 2500  # ARM!J122 (1750.6 -> 2500)
$ws.Cells.Item(122, 11).Value = 5796.8334  # ARM!K122 (5077.538399999999 -> 5796.8334)
$ws.Cells.Item(122, 12).Value = 7500  # ARM!L122 (5251.799999999999 -> 7500)
$ws.Cells.Item(122, 13).Value = -3346.8334  # ARM!M122 (-2627.538399999999 -> -3346.8334)
$ws.Cells.Item(122, 14).Value = -12400  # ARM!N122 (-10151.8 -> -12400)

$ws.Cells.Item(132, 8).Value = 2789.52  # ARM!H132 (1238.2354 -> 2789.52)
$ws.Cells.Item(132, 9).Value = 2558.9524  # ARM!I132 (1047.0333 -> 2558.9524)
$ws.Cells.Item(132, 10).Value = 4000  # ARM!J132 (2672.25 -> 4000)
$ws.Cells.Item(132, 11).Value = 7676.8572  # ARM!K132 (3141.0999 -> 7676.8572)
$ws.Cells.Item(132, 12).Value = 12000  # ARM!L132 (8016.75 -> 12000)
$ws.Cells.Item(132, 13).Value = -5146.8572  # ARM!M132 (-611.0999000000002 -> -5146.8572)
$ws.Cells.Item(132, 14).Value = -17060  # ARM!N132 (-13076.75 -> -17060)

$ws.Cells.Item(136, 8).Value = 6713.0415  # ARM!H136 (4456.4873 -> 6713.0415)
$ws.Cells.Item(136, 9).Value = 1711.2106  # ARM!I136 (1203.1212 -> 1711.2106)
$ws.Cells.Item(136, 10).Value = 25720  # ARM!J136 (22350 -> 25720)
$ws.Cells.Item(136, 11).Value = 5133.6318  # ARM!K136 (3609.3636 -> 5133.6318)
$ws.Cells.Item(136, 12).Value = 77160  # ARM!L136 (67050 -> 77160)
$ws.Cells.Item(136, 13).Value = -2583.6318  # ARM!M136 (-1059.3636 -> -2583.6318)
$ws.Cells.Item(136, 14).Value = -82260  # ARM!N136 (-72150 -> -82260)

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1134.7949  # BSM!H86 (1159.421 -> 1134.7949)
$ws.Cells.Item(86, 9).Value = 1148.4073  # BSM!I86 (1184.9231 -> 1148.4073)
$ws.Cells.Item(86, 11).Value = 1148.4073  # BSM!K86 (1184.9231 -> 1148.4073)
$ws.Cells.Item(86, 13).Value = -25.40730000000008  # BSM!M86 (-61.92309999999998 -> -25.40730000000008)

$ws.Cells.Item(89, 8).Value = 1134.7949  # BSM!H89 (1159.421 -> 1134.7949)
$ws.Cells.Item(89, 9).Value = 1148.4073  # BSM!I89 (1184.9231 -> 1148.4073)
$ws.Cells.Item(89, 11).Value = 5742.0365  # BSM!K89 (5924.6155 -> 5742.0365)
$ws.Cells.Item(89, 13).Value = -126.0365000000002  # BSM!M89 (-308.6154999999999 -> -126.0365000000002)

$ws.Cells.Item(105, 8).Value = 12491.5625  # BSM!H105 (11670.889 -> 12491.5625)
$ws.Cells.Item(105, 9).Value = 13583.214  # BSM!I105 (12317.25 -> 13583.214)
$ws.Cells.Item(105, 10).Value = 4850  # BSM!J105 (6500 -> 4850)
$ws.Cells.Item(105, 11).Value = 13583.214  # BSM!K105 (12317.25 -> 13583.214)
$ws.Cells.Item(105, 12).Value = 4850  # BSM!L105 (6500 -> 4850)
$ws.Cells.Item(105, 13).Value = -11836.214  # BSM!M105 (-10570.25 -> -11836.214)
$ws.Cells.Item(105, 14).Value = -8344  # BSM!N105 (-9994 -> -8344)

$ws.Cells.Item(134, 8).Value = 3033.4849  # BSM!H134 (3992.9333 -> 3033.4849)
$ws.Cells.Item(134, 9).Value = 2825.4644  # BSM!I134 (3908.0833 -> 2825.4644)
$ws.Cells.Item(134, 10).Value = 4198.4  # BSM!J134 (4332.3335 -> 4198.4)
$ws.Cells.Item(134, 11).Value = 8476.393199999999  # BSM!K134 (11724.2499 -> 8476.393199999999)
$ws.Cells.Item(134, 12).Value = 12595.2  # BSM!L134 (12997.0005 -> 12595.2)
$ws.Cells.Item(134, 13).Value = -5941.393199999999  # BSM!M134 (-9189.249899999999 -> -5941.393199999999)
$ws.Cells.Item(134, 14).Value = -17665.2  # BSM!N134 (-18067.0005 -> -17665.2)

$ws.Cells.Item(140, 8).Value = 101340.664  # BSM!H140 (100893.89 -> 101340.664)
$ws.Cells.Item(140, 10).Value = 101340.664  # BSM!J140 (100893.89 -> 101340.664)
$ws.Cells.Item(140, 12).Value = 101340.664  # BSM!L140 (100893.89 -> 101340.664)
$ws.Cells.Item(140, 14).Value = -111700.664  # BSM!N140 (-111253.89 -> -111700.664)

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 798.6667  # CRP!H58 (1061.9375 -> 798.6667)
$ws.Cells.Item(58, 9).Value = 780.5405  # CRP!I58 (1059.4 -> 780.5405)
$ws.Cells.Item(58, 10).Value = 932.8  # CRP!J58 (1100 -> 932.8)
$ws.Cells.Item(58, 11).Value = 780.5405  # CRP!K58 (1059.4 -> 780.5405)
$ws.Cells.Item(58, 12).Value = 932.8  # CRP!L58 (1100 -> 932.8)
$ws.Cells.Item(58, 13).Value = -577.5405  # CRP!M58 (-856.4000000000001 -> -577.5405)
$ws.Cells.Item(58, 14).Value = -1338.8  # CRP!N58 (-1506 -> -1338.8)

$ws.Cells.Item(99, 8).Value = 5179.2  # CRP!H99 (7965 -> 5179.2)
$ws.Cells.Item(99, 9).Value = 4768.8125  # CRP!I99 (7361.125 -> 4768.8125)
$ws.Cells.Item(99, 10).Value = 6820.75  # CRP!J99 (8770.166999999999 -> 6820.75)
$ws.Cells.Item(99, 11).Value = 4768.8125  # CRP!K99 (7361.125 -> 4768.8125)
$ws.Cells.Item(99, 12).Value = 6820.75  # CRP!L99 (8770.166999999999 -> 6820.75)
$ws.Cells.Item(99, 13).Value = -3270.8125  # CRP!M99 (-5863.125 -> -3270.8125)
$ws.Cells.Item(99, 14).Value = -9816.75  # CRP!N99 (-11766.167 -> -9816.75)

$ws.Cells.Item(126, 8).Value = 5179.2  # CRP!H126 (7965 -> 5179.2)
$ws.Cells.Item(126, 9).Value = 4768.8125  # CRP!I126 (7361.125 -> 4768.8125)
$ws.Cells.Item(126, 10).Value = 6820.75  # CRP!J126 (8770.166999999999 -> 6820.75)
$ws.Cells.Item(126, 11).Value = 14306.4375  # CRP!K126 (22083.375 -> 14306.4375)
$ws.Cells.Item(126, 12).Value = 20462.25  # CRP!L126 (26310.501 -> 20462.25)
$ws.Cells.Item(126, 13).Value = -11836.4375  # CRP!M126 (-19613.375 -> -11836.4375)
$ws.Cells.Item(126, 14).Value = -25402.25  # CRP!N126 (-31250.501 -> -25402.25)

$ws.Cells.Item(132, 8).Value = 24849.924  # CRP!H132 (45468.355 -> 24849.924)
$ws.Cells.Item(132, 9).Value = 36671.176  # CRP!I132 (87950.86 -> 36671.176)
$ws.Cells.Item(132, 10).Value = 2520.889  # CRP!J132 (2985.8572 -> 2520.889)
$ws.Cells.Item(132, 11).Value = 110013.528  # CRP!K132 (263852.58 -> 110013.528)
$ws.Cells.Item(132, 12).Value = 7562.667  # CRP!L132 (8957.571599999999 -> 7562.667)
$ws.Cells.Item(132, 13).Value = -107483.528  # CRP!M132 (-261322.58 -> -107483.528)
$ws.Cells.Item(132, 14).Value = -12622.667  # CRP!N132 (-14017.5716 -> -12622.667)

$ws.Cells.Item(134, 8).Value = 1508  # CRP!H134 (977.2593000000001 -> 1508)
$ws.Cells.Item(134, 9).Value = 1620  # CRP!I134 (934.3913 -> 1620)
$ws.Cells.Item(134, 10).Value = 500  # CRP!J134 (1223.75 -> 500)
$ws.Cells.Item(134, 11).Value = 4860  # CRP!K134 (2803.1739 -> 4860)
$ws.Cells.Item(134, 12).Value = 1500  # CRP!L134 (3671.25 -> 1500)
$ws.Cells.Item(134, 13).Value = -2325  # CRP!M134 (-268.1738999999998 -> -2325)
$ws.Cells.Item(134, 14).Value = -6570  # CRP!N134 (-8741.25 -> -6570)

$ws.Cells.Item(136, 8).Value = 798.6667  # CRP!H136 (1061.9375 -> 798.6667)
$ws.Cells.Item(136, 9).Value = 780.5405  # CRP!I136 (1059.4 -> 780.5405)
$ws.Cells.Item(136, 10).Value = 932.8  # CRP!J136 (1100 -> 932.8)
$ws.Cells.Item(136, 11).Value = 2341.6215  # CRP!K136 (3178.2 -> 2341.6215)
$ws.Cells.Item(136, 12).Value = 2798.4  # CRP!L136 (3300 -> 2798.4)
$ws.Cells.Item(136, 13).Value = 208.3785000000003  # CRP!M136 (-628.2000000000003 -> 208.3785000000003)
$ws.Cells.Item(136, 14).Value = -7898.4  # CRP!N136 (-8400 -> -7898.4)

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 71261816  # CUL!H4 (48973704 -> 71261816)
$ws.Cells.Item(4, 9).Value = 49381400  # CUL!I4 (65083652 -> 49381400)
$ws.Cells.Item(4, 10).Value = 121586780  # CUL!J4 (21586792 -> 121586780)
$ws.Cells.Item(4, 11).Value = 148144200  # CUL!K4 (195250956 -> 148144200)
$ws.Cells.Item(4, 12).Value = 364760340  # CUL!L4 (64760376 -> 364760340)
$ws.Cells.Item(4, 13).Value = -148144088  # CUL!M4 (-195250844 -> -148144088)
$ws.Cells.Item(4, 14).Value = -364760564  # CUL!N4 (-64760600 -> -364760564)

$ws.Cells.Item(12, 8).Value = 174.28572  # CUL!H12 (66.333336 -> 174.28572)
$ws.Cells.Item(12, 9).Value = 290.5  # CUL!I12 (40.75 -> 290.5)
$ws.Cells.Item(12, 10).Value = 127.8  # CUL!J12 (73.64286 -> 127.8)
$ws.Cells.Item(12, 11).Value = 871.5  # CUL!K12 (122.25 -> 871.5)
$ws.Cells.Item(12, 12).Value = 383.4  # CUL!L12 (220.92858 -> 383.4)
$ws.Cells.Item(12, 13).Value = -698.5  # CUL!M12 (50.75 -> -698.5)
$ws.Cells.Item(12, 14).Value = -729.4  # CUL!N12 (-566.92858 -> -729.4)

$ws.Cells.Item(132, 8).Value = 1806.2609  # CUL!H132 (1838.3636 -> 1806.2609)
$ws.Cells.Item(132, 10).Value = 1507.8235  # CUL!J132 (1533.3125 -> 1507.8235)
$ws.Cells.Item(132, 12).Value = 13570.4115  # CUL!L132 (13799.8125 -> 13570.4115)
$ws.Cells.Item(132, 14).Value = -18630.4115  # CUL!N132 (-18859.8125 -> -18630.4115)

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4237.4136  # GSM!H70 (4456.8184 -> 4237.4136)
$ws.Cells.Item(70, 9).Value = 4381.7915  # GSM!I70 (4548.1816 -> 4381.7915)
$ws.Cells.Item(70, 10).Value = 3544.4  # GSM!J70 (4274.091 -> 3544.4)
$ws.Cells.Item(70, 11).Value = 4381.7915  # GSM!K70 (4548.1816 -> 4381.7915)
$ws.Cells.Item(70, 12).Value = 3544.4  # GSM!L70 (4274.091 -> 3544.4)
$ws.Cells.Item(70, 13).Value = -4111.7915  # GSM!M70 (-4278.1816 -> -4111.7915)
$ws.Cells.Item(70, 14).Value = -4084.4  # GSM!N70 (-4814.091 -> -4084.4)

$ws.Cells.Item(73, 8).Value = 4237.4136  # GSM!H73 (4456.8184 -> 4237.4136)
$ws.Cells.Item(73, 9).Value = 4381.7915  # GSM!I73 (4548.1816 -> 4381.7915)
$ws.Cells.Item(73, 10).Value = 3544.4  # GSM!J73 (4274.091 -> 3544.4)
$ws.Cells.Item(73, 11).Value = 4381.7915  # GSM!K73 (4548.1816 -> 4381.7915)
$ws.Cells.Item(73, 12).Value = 3544.4  # GSM!L73 (4274.091 -> 3544.4)
$ws.Cells.Item(73, 13).Value = -3445.7915  # GSM!M73 (-3612.1816 -> -3445.7915)
$ws.Cells.Item(73, 14).Value = -5416.4  # GSM!N73 (-6146.091 -> -5416.4)

$ws.Cells.Item(102, 8).Value = 29764.555  # GSM!H102 (17583.227 -> 29764.555)
$ws.Cells.Item(102, 9).Value = 29764.555  # GSM!I102 (18136 -> 29764.555)
$ws.Cells.Item(102, 10).Value = 0  # GSM!J102 (1000 -> 0)
$ws.Cells.Item(102, 11).Value = 29764.555  # GSM!K102 (18136 -> 29764.555)
$ws.Cells.Item(102, 12).Value = 0  # GSM!L102 (1000 -> 0)
$ws.Cells.Item(102, 13).Value = -28142.555  # GSM!M102 (-16514 -> -28142.555)
$ws.Cells.Item(102, 14).ClearContents()  # GSM!N102 removed (was -4244)

$ws.Cells.Item(122, 8).Value = 1944.3948  # GSM!H122 (2046.2778 -> 1944.3948)
$ws.Cells.Item(122, 9).Value = 1821.6389  # GSM!I122 (1922.2941 -> 1821.6389)
$ws.Cells.Item(122, 11).Value = 5464.9167  # GSM!K122 (5766.8823 -> 5464.9167)
$ws.Cells.Item(122, 13).Value = -3014.9167  # GSM!M122 (-3316.8823 -> -3014.9167)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 2390  # LTW!H136 (3747.074 -> 2390)
$ws.Cells.Item(136, 9).Value = 2310  # LTW!I136 (3616.84 -> 2310)
$ws.Cells.Item(136, 10).Value = 5750  # LTW!J136 (5375 -> 5750)
$ws.Cells.Item(136, 11).Value = 6930  # LTW!K136 (10850.52 -> 6930)
$ws.Cells.Item(136, 12).Value = 17250  # LTW!L136 (16125 -> 17250)
$ws.Cells.Item(136, 13).Value = -4380  # LTW!M136 (-8300.52 -> -4380)
$ws.Cells.Item(136, 14).Value = -22350  # LTW!N136 (-21225 -> -22350)

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 16009  # WVR!H14 (4000 -> 16009)
$ws.Cells.Item(14, 9).Value = 0  # WVR!I14 (4000 -> 0)
$ws.Cells.Item(14, 10).Value = 16009  # WVR!J14 (0 -> 16009)
$ws.Cells.Item(14, 11).Value = 0  # WVR!K14 (4000 -> 0)
$ws.Cells.Item(14, 12).Value = 16009  # WVR!L14 (0 -> 16009)
$ws.Cells.Item(14, 13).ClearContents()  # WVR!M14 removed (was -3832)
$ws.Cells.Item(14, 14).Value = -16345  # WVR!N14 (None -> -16345)

$ws.Cells.Item(132, 8).Value = 5574587.5  # WVR!H132 (32950.27 -> 5574587.5)
$ws.Cells.Item(132, 9).Value = 6118080  # WVR!I132 (33868.28 -> 6118080)
$ws.Cells.Item(132, 10).Value = 3794.5  # WVR!J132 (10000 -> 3794.5)
$ws.Cells.Item(132, 11).Value = 18354240  # WVR!K132 (101604.84 -> 18354240)
$ws.Cells.Item(132, 12).Value = 11383.5  # WVR!L132 (30000 -> 11383.5)
$ws.Cells.Item(132, 13).Value = -18351710  # WVR!M132 (-99074.84 -> -18351710)
$ws.Cells.Item(132, 14).Value = -16443.5  # WVR!N132 (-35060 -> -16443.5)
